# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4288
$ws1.Range("F3").Value = 2432
$ws1.Range("F10").Value = 135
$ws1.Range("F12").Value = 1594
$ws1.Range("F14").Value = 3333
$ws1.Range("F15").Value = 225

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4288
$ws4.Range("F3").Value = 2432
$ws4.Range("F12").Value = 135
$ws4.Range("F16").Value = 1594
$ws4.Range("F18").Value = 3333
$ws4.Range("F19").Value = 225
